# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp text
$ws.Range("A1").Value = "Datos actualizados a 1 de Septiembre de 2020 a las 13:58"

# Row 61 - Suiza
$ws.Range("B61").Value = 42393
$ws.Range("C61").Value = 216
$ws.Range("E61").Value = 4283
$ws.Range("G61").Value = 4
$ws.Range("H61").Value = 2010

# Row 64 - Nepal
$ws.Range("B64").Value = 40529
$ws.Range("C64").Value = 1069
$ws.Range("D64").Value = 22178
$ws.Range("E64").Value = 18112
$ws.Range("G64").Value = 11
$ws.Range("H64").Value = 239

# Row 75 - Estado de Palestina
$ws.Range("B75").Value = 23281
$ws.Range("C75").Value = 552
$ws.Range("D75").Value = 15338
$ws.Range("E75").Value = 7784
$ws.Range("G75").Value = 7
$ws.Range("H75").Value = 159

# Row 84 - Madagascar
$ws.Range("B84").Value = 14957
$ws.Range("C84").Value = 94
$ws.Range("D84").Value = 13915
$ws.Range("E84").Value = 847
$ws.Range("G84").Value = 3
$ws.Range("H84").Value = 195

# Row 102 - Maldivas
$ws.Range("E102").Value = 2620
$ws.Range("G102").Value = 1
$ws.Range("H102").Value = 29

# Row 184 - Gibraltar
$ws.Range("B184").Value = 290
$ws.Range("C184").Value = 2
$ws.Range("D184").Value = 237
$ws.Range("E184").Value = 53
